$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 42.923077
$ws.Range("I9").Value = 36
$ws.Range("J9").Value = 81
$ws.Range("K9").Value = 36
$ws.Range("L9").Value = 81
$ws.Range("M9").Value = 133
$ws.Range("N9").Value = -419

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 499.625
$ws.Range("I33").Value = 259
$ws.Range("K33").Value = 259
$ws.Range("M33").Value = -30

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2251
$ws.Range("I76").Value = 2251
$ws.Range("K76").Value = 2251
$ws.Range("M76").Value = -1936

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 2251
$ws.Range("I79").Value = 2251
$ws.Range("K79").Value = 2251
$ws.Range("M79").Value = -1159

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1394.25
$ws.Range("I92").Value = 1394.25
$ws.Range("K92").Value = 1394.25
$ws.Range("M92").Value = -146.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 774.75
$ws.Range("I111").Value = 816.3333
$ws.Range("J111").Value = 650
$ws.Range("K111").Value = 2448.9999
$ws.Range("L111").Value = 1950
$ws.Range("M111").Value = 618.0001000000002
$ws.Range("N111").Value = -8084

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1833.25
$ws.Range("I127").Value = 1363.5454
$ws.Range("K127").Value = 4090.6362
$ws.Range("M127").Value = 869.3638000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1194.0476
$ws.Range("J138").Value = 1999
$ws.Range("L138").Value = 5997
$ws.Range("N138").Value = -16277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 30009
$ws.Range("J9").Value = 30009
$ws.Range("L9").Value = 30009
$ws.Range("N9").Value = -30349

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 550
$ws.Range("I19").Value = 550
$ws.Range("K19").Value = 550
$ws.Range("M19").Value = -321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 30009
$ws.Range("J20").Value = 30009
$ws.Range("L20").Value = 30009
$ws.Range("N20").Value = -30549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9817.866
$ws.Range("I32").Value = 8376.286
$ws.Range("K32").Value = 8376.286
$ws.Range("M32").Value = -8089.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2382.2856
$ws.Range("I122").Value = 2382.2856
$ws.Range("K122").Value = 7146.8568
$ws.Range("M122").Value = -4696.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5464.8
$ws.Range("I20").Value = 4331
$ws.Range("K20").Value = 4331
$ws.Range("M20").Value = -4084

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2490.5
$ws.Range("I86").Value = 2613.75
$ws.Range("K86").Value = 2613.75
$ws.Range("M86").Value = -1490.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2490.5
$ws.Range("I89").Value = 2613.75
$ws.Range("K89").Value = 13068.75
$ws.Range("M89").Value = -7452.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 100000
$ws.Range("J110").Value = 100000
$ws.Range("L110").Value = 100000
$ws.Range("N110").Value = -108180

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 39165.332
$ws.Range("J135").Value = 39165.332
$ws.Range("L135").Value = 39165.332
$ws.Range("N135").Value = -49305.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 688.1429000000001
$ws.Range("I19").Value = 376.52942
$ws.Range("J19").Value = 2012.5
$ws.Range("K19").Value = 376.52942
$ws.Range("L19").Value = 2012.5
$ws.Range("M19").Value = -206.52942
$ws.Range("N19").Value = -2352.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 688.1429000000001
$ws.Range("I24").Value = 376.52942
$ws.Range("J24").Value = 2012.5
$ws.Range("K24").Value = 376.52942
$ws.Range("L24").Value = 2012.5
$ws.Range("M24").Value = -206.52942
$ws.Range("N24").Value = -2352.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2013.8
$ws.Range("J94").Value = 1155
$ws.Range("L94").Value = 1155
$ws.Range("N94").Value = -2057

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 941.5833
$ws.Range("I107").Value = 982.2
$ws.Range("K107").Value = 982.2
$ws.Range("M107").Value = 937.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1157.3334
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 750
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 750
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 2250
$ws.Range("N41").Value = -2926
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 677.8
$ws.Range("I113").Value = 190
$ws.Range("J113").Value = 799.75
$ws.Range("K113").Value = 570
$ws.Range("L113").Value = 2399.25
$ws.Range("M113").Value = 1600
$ws.Range("N113").Value = -6739.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3998
$ws.Range("I139").Value = 3998
$ws.Range("K139").Value = 11994
$ws.Range("M139").Value = -6854

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1638
$ws.Range("I140").Value = 1638
$ws.Range("K140").Value = 4914
$ws.Range("M140").Value = 266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1954
$ws.Range("J10").Value = 1954
$ws.Range("L10").Value = 1954
$ws.Range("N10").Value = -2292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5462819
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 5462819
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 5462819
$ws.Range("N11").Value = -5463097
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1650.1428
$ws.Range("I126").Value = 1425.1666
$ws.Range("K126").Value = 4275.4998
$ws.Range("M126").Value = -1805.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6133.8887
$ws.Range("I40").Value = 3041
$ws.Range("K40").Value = 3041
$ws.Range("M40").Value = -2905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4168
$ws.Range("I122").Value = 3504
$ws.Range("K122").Value = 10512
$ws.Range("M122").Value = -8062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2395.5715
$ws.Range("I126").Value = 2104.6
$ws.Range("K126").Value = 6313.799999999999
$ws.Range("M126").Value = -3843.799999999999
